$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 185, shifting existing rows 185-233 down to 187-235.
$ws.Rows.Item(185).Resize(2).Insert()

# Fill in the new row 185 (Super Queen, Primera)
$ws.Cells.Item(185, 1).Value = 7
$ws.Cells.Item(185, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(185, 3).Value = "Ñuble"
$ws.Cells.Item(185, 4).Value = 44551
$ws.Cells.Item(185, 5).Value = 16
$ws.Cells.Item(185, 6).Value = "Fruta"
$ws.Cells.Item(185, 7).Value = 100103
$ws.Cells.Item(185, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(185, 9).Value = 100103006
$ws.Cells.Item(185, 10).Value = "Nectarín"
$ws.Cells.Item(185, 11).Value = "Super Queen"
$ws.Cells.Item(185, 12).Value = "Primera"
$ws.Cells.Item(185, 13).Value = 160
$ws.Cells.Item(185, 14).Value = 13000
$ws.Cells.Item(185, 15).Value = 14000
$ws.Cells.Item(185, 16).Value = 13500
$ws.Cells.Item(185, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(185, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(185, 19).Value = 844
$ws.Cells.Item(185, 20).Value = 16

# Fill in the new row 186 (Super Queen, Segunda)
$ws.Cells.Item(186, 1).Value = 7
$ws.Cells.Item(186, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(186, 3).Value = "Ñuble"
$ws.Cells.Item(186, 4).Value = 44551
$ws.Cells.Item(186, 5).Value = 16
$ws.Cells.Item(186, 6).Value = "Fruta"
$ws.Cells.Item(186, 7).Value = 100103
$ws.Cells.Item(186, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(186, 9).Value = 100103006
$ws.Cells.Item(186, 10).Value = "Nectarín"
$ws.Cells.Item(186, 11).Value = "Super Queen"
$ws.Cells.Item(186, 12).Value = "Segunda"
$ws.Cells.Item(186, 13).Value = 120
$ws.Cells.Item(186, 14).Value = 11000
$ws.Cells.Item(186, 15).Value = 12000
$ws.Cells.Item(186, 16).Value = 11500
$ws.Cells.Item(186, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(186, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(186, 19).Value = 719
$ws.Cells.Item(186, 20).Value = 16

# Apply date style (with the date display format) to the D cells of the two new rows,
# copying from the row above.
$ws.Cells.Item(184, 4).Copy()
$ws.Range($ws.Cells.Item(185, 4), $ws.Cells.Item(186, 4)).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
